$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9963980317115784
$ws.Range("B1").Value = 1.878901481628418
$ws.Range("C1").Value = 2.097723722457886
$ws.Range("D1").Value = 2.130523920059204
$ws.Range("E1").Value = 1.396522164344788
